$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch the existing wrap-text style (already used by F169) so that the
# new cells below which need the same "wrapped number" formatting dedupe
# onto that same style slot instead of Excel minting a brand new one.
$ws.Cells.Item(169, 6).WrapText = $true

# --- Row 170: tracker_date 25.03.2024 -------------------------------------
$ws.Cells.Item(170, 1).Value = "25.03.2024"
$ws.Cells.Item(170, 2).Value = "25.03.2024"
$ws.Cells.Item(170, 3).Value = 32333
$ws.Cells.Item(170, 3).WrapText = $true
$ws.Cells.Item(170, 4).Value = 13000
$ws.Cells.Item(170, 5).Value = 8400
$ws.Cells.Item(170, 6).Value = 74694
$ws.Cells.Item(170, 6).WrapText = $true
$ws.Cells.Item(170, 7).Value = 8663
$ws.Cells.Item(170, 8).Value = 6327
$ws.Cells.Item(170, 9).Value = 8000
$ws.Cells.Item(170, 10).Value = 450
$ws.Cells.Item(170, 11).Value = 116
$ws.Cells.Item(170, 12).Value = 4700
$ws.Cells.Item(170, 13).Value = "https://web.archive.org/web/20240325135227/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# --- Row 171: tracker_date 26.03.2024 (report_date still 25.03.2024) -----
$ws.Cells.Item(171, 1).Value = "26.03.2024"
$ws.Cells.Item(171, 2).Value = "25.03.2024"
$ws.Cells.Item(171, 3).Value = 32333
$ws.Cells.Item(171, 3).WrapText = $true
$ws.Cells.Item(171, 4).Value = 13000
$ws.Cells.Item(171, 5).Value = 8400
$ws.Cells.Item(171, 6).Value = 74694
$ws.Cells.Item(171, 6).WrapText = $true
$ws.Cells.Item(171, 7).Value = 8663
$ws.Cells.Item(171, 8).Value = 6327
$ws.Cells.Item(171, 9).Value = 8000
$ws.Cells.Item(171, 10).Value = 450
$ws.Cells.Item(171, 11).Value = 116
$ws.Cells.Item(171, 12).Value = 4700
$ws.Cells.Item(171, 13).Value = "https://web.archive.org/web/20240326214610/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# --- View bookkeeping: keep the frozen/split pane + selection pointed at
# the new bottom of the table, matching how the sheet looked after the
# rows were appended.
$excel.ActiveWindow.SplitRow = 155
$ws.Range("M171").Select() | Out-Null
